# Add a new row to the "Error" table ("表1") for a new error entry:
#   Id = 3001, Des = "资源不足" (Insufficient resources)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing table (ListObject) by one row; this keeps the table
# range, autoFilter range, and sheet dimension in sync automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Fill in the values for the newly added table row (row 24).
$ws.Cells.Item(24, 1).Value = 3001
$ws.Cells.Item(24, 2).Value = "资源不足"

# Match the author's final selection in the sheet.
$ws.Range("B22").Select()
